$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Clarence Kuhn
$ws.Range("A6").Value = 'Clarence'
$ws.Range("B6").Value = 'Kuhn'
$ws.Range("C6").Value = 'vanita.romaguera@hotmail.com'
$ws.Range("D6").Value = "'5097499366"
$ws.Range("E6").Value = 'Scientist'
$ws.Range("F6").Value = 't%iye7M%y'

# Row 7 - Jenifer VonRueden
$ws.Range("A7").Value = 'Jenifer'
$ws.Range("B7").Value = 'VonRueden'
$ws.Range("C7").Value = 'clark.harris@yahoo.com'
$ws.Range("D7").Value = "'9898922868"
$ws.Range("E7").Value = 'Engineer'
$ws.Range("F7").Value = 'oi$g#7Pi'

# Re-apply the plain (default) cell format used by the rest of the table onto
# the two phone-number cells so they keep the same style as their siblings
# instead of picking up the "quote prefix" formatting that typing a leading
# apostrophe applies.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
